$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Bundesliga row values from 1773 to 1761
$ws.Range("B5").Value = 1761
$ws.Range("C5").Value = 1761

# Update the selected cell/range on the sheet
$ws.Range("C6").Select()
